$d = $word.ActiveDocument

function Get-ParaByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs($i)
        if ($p.Range.Text -eq $text) {
            return $p
        }
    }
    return $null
}

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

# ---------------------------------------------------------------------------
# STEP A: "Screenshake" paragraph gets a <w:lastRenderedPageBreak/> before its text
# ---------------------------------------------------------------------------
$p = Get-ParaByText $d "Screenshake`r"
$body = '<w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t>Screenshake</w:t></w:r></w:p></w:body>'
$p.Range.InsertXML($pkgHeader + $body + $pkgFooter)

# ---------------------------------------------------------------------------
# STEP B: insert the "Modular enemy" block right before the "Feel -" paragraph
# ---------------------------------------------------------------------------
$p = Get-ParaByText $d "Feel –`r"
$p.Range.InsertParagraphBefore()
$insertedPara = Get-ParaByText $d "Feel –`r"
$target = $insertedPara.Previous().Range
$body = '<w:body>' + `
  '<w:p><w:r><w:t xml:space="preserve">Modular enemy – </w:t></w:r></w:p>' + `
  '<w:p><w:r><w:t>Main piece – Spawns modules</w:t></w:r><w:r><w:t>, attaches modules to itself, can move and has AI</w:t></w:r></w:p>' + `
  '<w:p><w:r><w:tab/><w:t xml:space="preserve">Module types – </w:t></w:r></w:p>' + `
  '<w:p><w:r><w:tab/></w:r><w:r><w:tab/><w:t>Module holder – spawns more modules</w:t></w:r></w:p>' + `
  '<w:p><w:r><w:tab/></w:r><w:r><w:tab/><w:t>Armour – Protects the core</w:t></w:r></w:p>' + `
  '<w:p><w:r><w:tab/></w:r><w:r><w:tab/><w:t>Weak spot – Explodes and damages other parts</w:t></w:r></w:p>' + `
  '<w:p><w:r><w:tab/></w:r><w:r><w:tab/><w:t>Machine gun – Rapid fire, weak bullets</w:t></w:r></w:p>' + `
  '<w:p><w:r><w:tab/></w:r><w:r><w:tab/><w:t>Missiles – Large, slow projectiles that explode and can remove cover</w:t></w:r></w:p>' + `
  '<w:p/>' + `
  '</w:body>'
$target.InsertXML($pkgHeader + $body + $pkgFooter)

# ---------------------------------------------------------------------------
# STEP C: rebuild "Scan...", blank, "Movement - " (+bookmark) + "Move, roll..."
# ---------------------------------------------------------------------------
$p = Get-ParaByText $d "Move, roll, jump, slide, climb`r"
$body = '<w:body>' + `
  '<w:p><w:r><w:t xml:space="preserve">Movement - </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>Move, roll, jump, slide, climb</w:t></w:r></w:p>' + `
  '</w:body>'
$p.Range.InsertXML($pkgHeader + $body + $pkgFooter)

# ---------------------------------------------------------------------------
# STEP D: "traps" becomes two runs "T" + "raps" (same rendered text)
# ---------------------------------------------------------------------------
$p = Get-ParaByText $d "traps`r"
$body = '<w:body><w:p><w:r><w:t>T</w:t></w:r><w:r><w:t>raps</w:t></w:r></w:p></w:body>'
$p.Range.InsertXML($pkgHeader + $body + $pkgFooter)

# ---------------------------------------------------------------------------
# STEP E: strip the _GoBack bookmark off the end of the "Projectile types" paragraph
# ---------------------------------------------------------------------------
$p = Get-ParaByText $d "Projectile types – Tear, fire DOT, corruption Hack, shock Paralysis, freeze increased damage`r"
$body = '<w:body><w:p>' + `
  '<w:r><w:t xml:space="preserve">Projectile types – </w:t></w:r>' + `
  '<w:r><w:t>Tear, fire</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> DOT</w:t></w:r>' + `
  '<w:r><w:t>, corruption</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> Hack</w:t></w:r>' + `
  '<w:r><w:t>, shock</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> Paralysis</w:t></w:r>' + `
  '<w:r><w:t>, freeze</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> increased damage</w:t></w:r>' + `
  '</w:p></w:body>'
$p.Range.InsertXML($pkgHeader + $body + $pkgFooter)

# ---------------------------------------------------------------------------
# STEP F: new first paragraph "Notes on HZD" before "Character"
# ---------------------------------------------------------------------------
$p = Get-ParaByText $d "Character`r"
$p.Range.InsertParagraphBefore()
(Get-ParaByText $d "Character`r").Previous().Range.Text = "Notes on HZD"
